$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H97").Value = 4679.25
$ws.Range("J97").Value = 5319.143
$ws.Range("L97").Value = 15957.429
$ws.Range("N97").Value = -16949.429
$ws.Range("H106").Value = 3332.5
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 3332.5
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 3332.5
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -4594.5
$ws.Range("H132").Value = 1514.1666
$ws.Range("I132").Value = 1328.8125
$ws.Range("J132").Value = 2997
$ws.Range("K132").Value = 3986.4375
$ws.Range("L132").Value = 8991
$ws.Range("M132").Value = -1456.4375
$ws.Range("N132").Value = -14051
$ws.Range("H138").Value = 4496.4
$ws.Range("J138").Value = 3600.3696
$ws.Range("L138").Value = 10801.1088
$ws.Range("N138").Value = -21081.1088

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 722.5357
$ws.Range("I2").Value = 750.38464
$ws.Range("K2").Value = 750.38464
$ws.Range("M2").Value = -637.38464
$ws.Range("H14").Value = 12763
$ws.Range("I14").Value = 423
$ws.Range("J14").Value = 17699
$ws.Range("K14").Value = 423
$ws.Range("L14").Value = 17699
$ws.Range("M14").Value = -248
$ws.Range("N14").Value = -18049
$ws.Range("H25").Value = 18977
$ws.Range("I25").Value = 5899.5
$ws.Range("J25").Value = 23336.166
$ws.Range("K25").Value = 5899.5
$ws.Range("L25").Value = 23336.166
$ws.Range("M25").Value = -5497.5
$ws.Range("N25").Value = -24140.166
$ws.Range("H61").Value = 786744.7
$ws.Range("I61").Value = 2133.2456
$ws.Range("J61").Value = 3271347.5
$ws.Range("K61").Value = 2133.2456
$ws.Range("L61").Value = 3271347.5
$ws.Range("M61").Value = -1921.2456
$ws.Range("N61").Value = -3271771.5
$ws.Range("H74").Value = 448327.25
$ws.Range("I74").Value = 1419.5625
$ws.Range("J74").Value = 1242829.9
$ws.Range("K74").Value = 1419.5625
$ws.Range("L74").Value = 1242829.9
$ws.Range("M74").Value = -545.5625
$ws.Range("N74").Value = -1244577.9
$ws.Range("H77").Value = 448327.25
$ws.Range("I77").Value = 1419.5625
$ws.Range("J77").Value = 1242829.9
$ws.Range("K77").Value = 7097.8125
$ws.Range("L77").Value = 6214149.5
$ws.Range("M77").Value = -2729.8125
$ws.Range("N77").Value = -6222885.5
$ws.Range("H116").Value = 722.5357
$ws.Range("I116").Value = 750.38464
$ws.Range("K116").Value = 750.38464
$ws.Range("M116").Value = 1543.61536
$ws.Range("H136").Value = 786744.7
$ws.Range("I136").Value = 2133.2456
$ws.Range("J136").Value = 3271347.5
$ws.Range("K136").Value = 6399.736800000001
$ws.Range("L136").Value = 9814042.5
$ws.Range("M136").Value = -3849.736800000001
$ws.Range("N136").Value = -9819142.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 722.5357
$ws.Range("I3").Value = 750.38464
$ws.Range("K3").Value = 750.38464
$ws.Range("M3").Value = -636.38464
$ws.Range("H80").Value = 83333570
$ws.Range("J80").Value = 367.7143
$ws.Range("L80").Value = 367.7143
$ws.Range("N80").Value = -2363.7143
$ws.Range("H83").Value = 83333570
$ws.Range("J83").Value = 367.7143
$ws.Range("L83").Value = 1838.5715
$ws.Range("N83").Value = -11822.5715
$ws.Range("H107").Value = 8472.728
$ws.Range("I107").Value = 9483.321
$ws.Range("K107").Value = 9483.321
$ws.Range("M107").Value = -7563.321

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2831.6345
$ws.Range("I31").Value = 2412.5806
$ws.Range("K31").Value = 2412.5806
$ws.Range("M31").Value = -2117.5806
$ws.Range("H34").Value = 2831.6345
$ws.Range("I34").Value = 2412.5806
$ws.Range("K34").Value = 2412.5806
$ws.Range("M34").Value = -2210.5806
$ws.Range("H58").Value = 1398.9454
$ws.Range("I58").Value = 1341.3846
$ws.Range("J58").Value = 1539.25
$ws.Range("K58").Value = 1341.3846
$ws.Range("L58").Value = 1539.25
$ws.Range("M58").Value = -1138.3846
$ws.Range("N58").Value = -1945.25
$ws.Range("H105").Value = 3137.8572
$ws.Range("I105").Value = 2158
$ws.Range("K105").Value = 2158
$ws.Range("M105").Value = -411
$ws.Range("H136").Value = 1398.9454
$ws.Range("I136").Value = 1341.3846
$ws.Range("J136").Value = 1539.25
$ws.Range("K136").Value = 4024.1538
$ws.Range("L136").Value = 4617.75
$ws.Range("M136").Value = -1474.1538
$ws.Range("N136").Value = -9717.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 6974.75
$ws.Range("J44").Value = 11950
$ws.Range("L44").Value = 35850
$ws.Range("N44").Value = -36646
$ws.Range("H113").Value = 470.64285
$ws.Range("I113").Value = 657.2
$ws.Range("K113").Value = 1971.6
$ws.Range("M113").Value = 198.3999999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H132").Value = 992709.9
$ws.Range("I132").Value = 7787.9443
$ws.Range("J132").Value = 3208784
$ws.Range("K132").Value = 23363.8329
$ws.Range("L132").Value = 9626352
$ws.Range("M132").Value = -20833.8329
$ws.Range("N132").Value = -9631412

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2515.4736
$ws.Range("I40").Value = 1349.9166
$ws.Range("J40").Value = 4513.5713
$ws.Range("K40").Value = 1349.9166
$ws.Range("L40").Value = 4513.5713
$ws.Range("M40").Value = -1213.9166
$ws.Range("N40").Value = -4785.5713
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 3110.4482
$ws.Range("I132").Value = 2765.95
$ws.Range("K132").Value = 8297.849999999999
$ws.Range("M132").Value = -5767.849999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 725.8
$ws.Range("I113").Value = 823.0833
$ws.Range("J113").Value = 579.875
$ws.Range("K113").Value = 2469.2499
$ws.Range("L113").Value = 1739.625
$ws.Range("M113").Value = -299.2498999999998
$ws.Range("N113").Value = -6079.625
$ws.Range("H126").Value = 3013.4666
$ws.Range("I126").Value = 2940.6155
$ws.Range("K126").Value = 8821.8465
$ws.Range("M126").Value = -6351.8465
$ws.Range("H141").Value = 299999
$ws.Range("J141").Value = 299999
$ws.Range("L141").Value = 299999
$ws.Range("N141").Value = -310359
